$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Volume/Number and report date range text (rich text shared strings) ---
$ws.Range("A8").Value = "Volume 29   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/14/2022  Through  11/20/2022"

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("G14").Value = 1
$ws.Range("M15").Value = 22.222222222222
$ws.Range("N15").Value = -59.259259259259
$ws.Range("C16").Value = 9
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 28.571428571428
$ws.Range("F16").Value = 50
$ws.Range("H16").Value = 13.636363636363
$ws.Range("I16").Value = 402
$ws.Range("J16").Value = 316
$ws.Range("K16").Value = 27.215189873417
$ws.Range("L16").Value = 54.615384615384
$ws.Range("M16").Value = 11.357340720221
$ws.Range("N16").Value = -61.084220716360
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -41.666666666666
$ws.Range("F17").Value = 45
$ws.Range("G17").Value = 44
$ws.Range("H17").Value = 2.272727272727
$ws.Range("I17").Value = 560
$ws.Range("J17").Value = 435
$ws.Range("K17").Value = 28.735632183908
$ws.Range("L17").Value = 39.650872817955
$ws.Range("M17").Value = 65.191740412979
$ws.Range("N17").Value = -10.4
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -14.285714285714
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 237
$ws.Range("J18").Value = 198
$ws.Range("K18").Value = 19.696969696969
$ws.Range("L18").Value = 7.239819004524
$ws.Range("M18").Value = 5.333333333333
$ws.Range("N18").Value = -79.795396419437
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -9.090909090909
$ws.Range("F19").Value = 46
$ws.Range("H19").Value = 31.428571428571
$ws.Range("I19").Value = 418
$ws.Range("J19").Value = 413
$ws.Range("K19").Value = 1.210653753026
$ws.Range("L19").Value = 9.424083769633
$ws.Range("M19").Value = 45.138888888888
$ws.Range("N19").Value = 0.966183574879
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 66.666666666666
$ws.Range("F20").Value = 28
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = 75
$ws.Range("I20").Value = 263
$ws.Range("J20").Value = 173
$ws.Range("K20").Value = 52.023121387283
$ws.Range("L20").Value = 103.875968992248
$ws.Range("M20").Value = 141.284403669725
$ws.Range("N20").Value = -51.831501831501
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = -7.5
$ws.Range("F21").Value = 185
$ws.Range("G21").Value = 174
$ws.Range("H21").Value = 6.321839080459
$ws.Range("I21").Value = 1914
$ws.Range("J21").Value = 1581
$ws.Range("K21").Value = 21.062618595825
$ws.Range("L21").Value = 33.752620545073
$ws.Range("M21").Value = 42.410714285714
$ws.Range("N21").Value = -50.657385924207
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = -20
$ws.Range("H23").Value = -100
$ws.Range("L23").Value = -21.212121212121
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = 44.444444444444
$ws.Range("F24").Value = 84
$ws.Range("G24").Value = 83
$ws.Range("H24").Value = 1.204819277108
$ws.Range("I24").Value = 1099
$ws.Range("J24").Value = 798
$ws.Range("K24").Value = 37.719298245614
$ws.Range("L24").Value = 19.068255687974
$ws.Range("M24").Value = 26.321839080459
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = -10
$ws.Range("F25").Value = 56
$ws.Range("G25").Value = 60
$ws.Range("H25").Value = -6.666666666666
$ws.Range("I25").Value = 792
$ws.Range("J25").Value = 693
$ws.Range("K25").Value = 14.285714285714
$ws.Range("L25").Value = 8.641975308641
$ws.Range("M25").Value = -9.897610921501
$ws.Range("L26").Value = 9.756097560975
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -71.428571428571
$ws.Range("I27").Value = 65
$ws.Range("J27").Value = 50
$ws.Range("K27").Value = 30
$ws.Range("L27").Value = 14.035087719298
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 49
$ws.Range("K28").Value = -9.259259259259
$ws.Range("L28").Value = 6.521739130434
$ws.Range("M28").Value = 25.641025641025
$ws.Range("N28").Value = -60.483870967741
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -33.333333333333
$ws.Range("I29").Value = 43
$ws.Range("K29").Value = -12.244897959183
$ws.Range("L29").Value = 19.444444444444
$ws.Range("M29").Value = 34.375
$ws.Range("N29").Value = -58.653846153846

# --- Cells whose data type/style changes (number <-> text) ---
# These require re-applying the correct cell style (number format) after setting
# the value, since Excel infers a style purely from the assigned value/type.
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = 1
$ws.Range("C36").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = -100
$ws.Range("K36").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("D23").PasteSpecial(-4122)

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "***.*"
$ws.Range("A36").Copy()
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("F23").PasteSpecial(-4122)

$ws.Range("C27").Value = 1
$ws.Range("C36").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("C28").Value = 1
$ws.Range("C36").Copy()
$ws.Range("C28").PasteSpecial(-4122)

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("D28").PasteSpecial(-4122)

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("A36").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("C29").Value = 1
$ws.Range("C36").Copy()
$ws.Range("C29").PasteSpecial(-4122)

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("A36").Copy()
$ws.Range("D29").PasteSpecial(-4122)

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("A36").Copy()
$ws.Range("E29").PasteSpecial(-4122)
